$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.02817156464667
$ws.Cells.Item(2, 4).Value = 1.030020546426149
$ws.Cells.Item(2, 5).Value = 1.037310853597163
$ws.Cells.Item(2, 6).Value = 1.046259149266063
$ws.Cells.Item(2, 9).Value = 1.032317986453552
$ws.Cells.Item(2, 10).Value = 1.03332536267446
$ws.Cells.Item(2, 11).Value = 1.032832768017343
$ws.Cells.Item(2, 12).Value = 1.040102101284833
$ws.Cells.Item(2, 13).Value = 1.049025086651996
$ws.Cells.Item(2, 14).Value = 1.034792803213533

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.029407233558643
$ws.Cells.Item(3, 4).Value = 1.030461415834738
$ws.Cells.Item(3, 5).Value = 1.038418322210282
$ws.Cells.Item(3, 6).Value = 1.047482177024637
$ws.Cells.Item(3, 9).Value = 1.03241346095033
$ws.Cells.Item(3, 10).Value = 1.034200060393832
$ws.Cells.Item(3, 11).Value = 1.033082893208857
$ws.Cells.Item(3, 12).Value = 1.041018563157604
$ws.Cells.Item(3, 13).Value = 1.050058651655641
$ws.Cells.Item(3, 14).Value = 1.035668743103996

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.03020627175964
$ws.Cells.Item(4, 4).Value = 1.030744229107144
$ws.Cells.Item(4, 5).Value = 1.039134339492755
$ws.Cells.Item(4, 6).Value = 1.048271815459117
$ws.Cells.Item(4, 9).Value = 1.032471999206458
$ws.Cells.Item(4, 10).Value = 1.034765059320114
$ws.Cells.Item(4, 11).Value = 1.033241558178008
$ws.Cells.Item(4, 12).Value = 1.041610403523623
$ws.Cells.Item(4, 13).Value = 1.050725137056315
$ws.Cells.Item(4, 14).Value = 1.03623454439356

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030542064884434
$ws.Cells.Item(5, 4).Value = 1.030862532816558
$ws.Cells.Item(5, 5).Value = 1.039435213475144
$ws.Cells.Item(5, 6).Value = 1.048603363537542
$ws.Cells.Item(5, 9).Value = 1.032495831543465
$ws.Cells.Item(5, 10).Value = 1.035002349455757
$ws.Cells.Item(5, 11).Value = 1.033307497027641
$ws.Cells.Item(5, 12).Value = 1.04185893383425
$ws.Cells.Item(5, 13).Value = 1.051004777618448
$ws.Cells.Item(5, 14).Value = 1.036472171508413

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030598438913854
$ws.Cells.Item(6, 4).Value = 1.030882361803813
$ws.Cells.Item(6, 5).Value = 1.039485723387253
$ws.Cells.Item(6, 6).Value = 1.048659007562714
$ws.Cells.Item(6, 9).Value = 1.032499787482363
$ws.Cells.Item(6, 10).Value = 1.035042177747753
$ws.Cells.Item(6, 11).Value = 1.033318523597193
$ws.Cells.Item(6, 12).Value = 1.041900646816567
$ws.Cells.Item(6, 13).Value = 1.051051698274787
$ws.Cells.Item(6, 14).Value = 1.036512056361152

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030210759122265
$ws.Cells.Item(7, 4).Value = 1.030745812212391
$ws.Cells.Item(7, 5).Value = 1.039138360332957
$ws.Cells.Item(7, 6).Value = 1.048276247253745
$ws.Cells.Item(7, 9).Value = 1.032472320710888
$ws.Cells.Item(7, 10).Value = 1.03476823092593
$ws.Cells.Item(7, 11).Value = 1.033242442260326
$ws.Cells.Item(7, 12).Value = 1.041613725494255
$ws.Cells.Item(7, 13).Value = 1.050728875787133
$ws.Cells.Item(7, 14).Value = 1.03623772050342

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028589274370517
$ws.Cells.Item(8, 4).Value = 1.030170047995659
$ws.Cells.Item(8, 5).Value = 1.037685250550765
$ws.Cells.Item(8, 6).Value = 1.046672838712672
$ws.Cells.Item(8, 9).Value = 1.032350922536984
$ws.Cells.Item(8, 10).Value = 1.033621176462157
$ws.Cells.Item(8, 11).Value = 1.032917956856678
$ws.Cells.Item(8, 12).Value = 1.040412067178979
$ws.Cells.Item(8, 13).Value = 1.049374861560514
$ws.Cells.Item(8, 14).Value = 1.035089037090739

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.025727898959768
$ws.Cells.Item(9, 4).Value = 1.029136748421671
$ws.Cells.Item(9, 5).Value = 1.035120113790983
$ws.Cells.Item(9, 6).Value = 1.043834027412905
$ws.Cells.Item(9, 9).Value = 1.032112249757243
$ws.Cells.Item(9, 10).Value = 1.031592283380329
$ws.Cells.Item(9, 11).Value = 1.032321878765232
$ws.Cells.Item(9, 12).Value = 1.038285564544728
$ws.Cells.Item(9, 13).Value = 1.046971259594827
$ws.Cells.Item(9, 14).Value = 1.033057262748018

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.023817389417192
$ws.Cells.Item(10, 4).Value = 1.028435421405652
$ws.Cells.Item(10, 5).Value = 1.033406858807775
$ws.Cells.Item(10, 6).Value = 1.041932390233029
$ws.Cells.Item(10, 9).Value = 1.031936571484487
$ws.Cells.Item(10, 10).Value = 1.030234462923442
$ws.Cells.Item(10, 11).Value = 1.031908266848974
$ws.Cells.Item(10, 12).Value = 1.036861748026829
$ws.Cells.Item(10, 13).Value = 1.045356926122214
$ws.Cells.Item(10, 14).Value = 1.031697514030337

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.022989386706896
$ws.Cells.Item(11, 4).Value = 1.028128810055141
$ws.Cells.Item(11, 5).Value = 1.032664228554953
$ws.Cells.Item(11, 6).Value = 1.041106781368141
$ws.Cells.Item(11, 9).Value = 1.03185658813384
$ws.Cells.Item(11, 10).Value = 1.02964525102772
$ws.Cells.Item(11, 11).Value = 1.031725340812541
$ws.Cells.Item(11, 12).Value = 1.036243742582085
$ws.Cells.Item(11, 13).Value = 1.044655055907513
$ws.Cells.Item(11, 14).Value = 1.03110746538614

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022681715226343
$ws.Cells.Item(12, 4).Value = 1.028014482443979
$ws.Cells.Item(12, 5).Value = 1.032388263493682
$ws.Cells.Item(12, 6).Value = 1.040799782771414
$ws.Cells.Item(12, 9).Value = 1.03182629196342
$ws.Cells.Item(12, 10).Value = 1.029426199229555
$ws.Cells.Item(12, 11).Value = 1.031656820200518
$ws.Cells.Item(12, 12).Value = 1.036013963187344
$ws.Cells.Item(12, 13).Value = 1.044393919811891
$ws.Cells.Item(12, 14).Value = 1.030888102509296

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022747717033964
$ws.Cells.Item(13, 4).Value = 1.028039025916199
$ws.Cells.Item(13, 5).Value = 1.032447464359042
$ws.Cells.Item(13, 6).Value = 1.040865650011373
$ws.Cells.Item(13, 9).Value = 1.031832817130157
$ws.Cells.Item(13, 10).Value = 1.029473195316047
$ws.Cells.Item(13, 11).Value = 1.031671544042497
$ws.Cells.Item(13, 12).Value = 1.03606326184255
$ws.Cells.Item(13, 13).Value = 1.044449953886516
$ws.Cells.Item(13, 14).Value = 1.030935165335622

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.022963956839462
$ws.Cells.Item(14, 4).Value = 1.028119368633059
$ws.Cells.Item(14, 5).Value = 1.032641419638618
$ws.Cells.Item(14, 6).Value = 1.041081411526799
$ws.Cells.Item(14, 9).Value = 1.031854095806751
$ws.Cells.Item(14, 10).Value = 1.029627148071355
$ws.Cells.Item(14, 11).Value = 1.03171968857253
$ws.Cells.Item(14, 12).Value = 1.036224753538951
$ws.Cells.Item(14, 13).Value = 1.044633479102696
$ws.Cells.Item(14, 14).Value = 1.0310893367215

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.023097174060847
$ws.Cells.Item(15, 4).Value = 1.028168812385732
$ws.Cells.Item(15, 5).Value = 1.03276090604266
$ws.Cells.Item(15, 6).Value = 1.04121430543314
$ws.Cells.Item(15, 9).Value = 1.031867128580078
$ws.Cells.Item(15, 10).Value = 1.029721977908785
$ws.Cells.Item(15, 11).Value = 1.031749276026882
$ws.Cells.Item(15, 12).Value = 1.036324224080171
$ws.Cells.Item(15, 13).Value = 1.04474649798505
$ws.Cells.Item(15, 14).Value = 1.031184301228177

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.023872325368482
$ws.Cells.Item(16, 4).Value = 1.028455708588133
$ws.Cells.Item(16, 5).Value = 1.033456128137704
$ws.Cells.Item(16, 6).Value = 1.041987136912546
$ws.Cells.Item(16, 9).Value = 1.031941797363679
$ws.Cells.Item(16, 10).Value = 1.030273540122215
$ws.Cells.Item(16, 11).Value = 1.031920326503756
$ws.Cells.Item(16, 12).Value = 1.036902731606338
$ws.Cells.Item(16, 13).Value = 1.045403446689255
$ws.Cells.Item(16, 14).Value = 1.031736646723216

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.024358356999422
$ws.Cells.Item(17, 4).Value = 1.028634887425422
$ws.Cells.Item(17, 5).Value = 1.033892012482009
$ws.Cells.Item(17, 6).Value = 1.042471326603344
$ws.Cells.Item(17, 9).Value = 1.031987588201136
$ws.Cells.Item(17, 10).Value = 1.030619180024689
$ws.Cells.Item(17, 11).Value = 1.032026597857157
$ws.Cells.Item(17, 12).Value = 1.0372652156747
$ws.Cells.Item(17, 13).Value = 1.045814768273885
$ws.Cells.Item(17, 14).Value = 1.032082777473997

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.024641779904964
$ws.Cells.Item(18, 4).Value = 1.028739116371629
$ws.Cells.Item(18, 5).Value = 1.034146181369793
$ws.Cells.Item(18, 6).Value = 1.042753535348707
$ws.Cells.Item(18, 9).Value = 1.032013919550868
$ws.Cells.Item(18, 10).Value = 1.030820663851126
$ws.Cells.Item(18, 11).Value = 1.032088214561481
$ws.Cells.Item(18, 12).Value = 1.037476503486256
$ws.Cells.Item(18, 13).Value = 1.046054410017996
$ws.Cells.Item(18, 14).Value = 1.032284547430578

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.02473840782714
$ws.Cells.Item(19, 4).Value = 1.028774607677994
$ws.Cells.Item(19, 5).Value = 1.034232833717592
$ws.Cells.Item(19, 6).Value = 1.04284972552511
$ws.Cells.Item(19, 9).Value = 1.032022833764723
$ws.Cells.Item(19, 10).Value = 1.030889343980569
$ws.Cells.Item(19, 11).Value = 1.032109161529225
$ws.Cells.Item(19, 12).Value = 1.037548522965141
$ws.Cells.Item(19, 13).Value = 1.046136075031644
$ws.Cells.Item(19, 14).Value = 1.032353325093682

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02430621778689
$ws.Cells.Item(20, 4).Value = 1.028615692493618
$ws.Cells.Item(20, 5).Value = 1.033845253998067
$ws.Cells.Item(20, 6).Value = 1.042419399463508
$ws.Cells.Item(20, 9).Value = 1.031982714333703
$ws.Cells.Item(20, 10).Value = 1.030582108776806
$ws.Cells.Item(20, 11).Value = 1.032015234162721
$ws.Cells.Item(20, 12).Value = 1.037226339371019
$ws.Cells.Item(20, 13).Value = 1.045770665836141
$ws.Cells.Item(20, 14).Value = 1.03204565358069

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.022900282791307
$ws.Cells.Item(21, 4).Value = 1.028095721794414
$ws.Cells.Item(21, 5).Value = 1.032584307947534
$ws.Cells.Item(21, 6).Value = 1.041017884284153
$ws.Cells.Item(21, 9).Value = 1.031847845956053
$ws.Cells.Item(21, 10).Value = 1.029581818145445
$ws.Cells.Item(21, 11).Value = 1.031705527033981
$ws.Cells.Item(21, 12).Value = 1.036177204475474
$ws.Cells.Item(21, 13).Value = 1.044579447390736
$ws.Cells.Item(21, 14).Value = 1.031043942421896

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.022015652072926
$ws.Cells.Item(22, 4).Value = 1.027766259482751
$ws.Cells.Item(22, 5).Value = 1.031790810911937
$ws.Cells.Item(22, 6).Value = 1.040134781649689
$ws.Cells.Item(22, 9).Value = 1.031759654398088
$ws.Cells.Item(22, 10).Value = 1.028951782344199
$ws.Cells.Item(22, 11).Value = 1.031507483141244
$ws.Cells.Item(22, 12).Value = 1.035516271032875
$ws.Cells.Item(22, 13).Value = 1.043827990971906
$ws.Cells.Item(22, 14).Value = 1.030413011897547

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022484675568429
$ws.Cells.Item(23, 4).Value = 1.027941153424203
$ws.Cells.Item(23, 5).Value = 1.032211524805796
$ws.Cells.Item(23, 6).Value = 1.040603113235093
$ws.Cells.Item(23, 9).Value = 1.031806727788923
$ws.Cells.Item(23, 10).Value = 1.029285882441282
$ws.Cells.Item(23, 11).Value = 1.031612784021551
$ws.Cells.Item(23, 12).Value = 1.035866768249379
$ws.Cells.Item(23, 13).Value = 1.04422658881602
$ws.Cells.Item(23, 14).Value = 1.03074758645509

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.024329777472366
$ws.Cells.Item(24, 4).Value = 1.02862436673197
$ws.Cells.Item(24, 5).Value = 1.033866382376237
$ws.Cells.Item(24, 6).Value = 1.042442863754502
$ws.Cells.Item(24, 9).Value = 1.031984917791739
$ws.Cells.Item(24, 10).Value = 1.030598860055185
$ws.Cells.Item(24, 11).Value = 1.03202037006931
$ws.Cells.Item(24, 12).Value = 1.037243906341116
$ws.Cells.Item(24, 13).Value = 1.045790594679665
$ws.Cells.Item(24, 14).Value = 1.032062428647805

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.026468136884543
$ws.Cells.Item(25, 4).Value = 1.029406086815867
$ws.Cells.Item(25, 5).Value = 1.035783814041064
$ws.Cells.Item(25, 6).Value = 1.044569524920671
$ws.Cells.Item(25, 9).Value = 1.032317986453552
$ws.Cells.Item(25, 10).Value = 1.032117715078505
$ws.Cells.Item(25, 11).Value = 1.03247884686872
$ws.Cells.Item(25, 12).Value = 1.038836394193408
$ws.Cells.Item(25, 13).Value = 1.047594746742769
$ws.Cells.Item(25, 14).Value = 1.033583440619473
